$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends at row 29, which holds customer "09876543" (stored
# as text, with its leading zero) and 0 points. Insert a fresh row 30 below it
# so the original text record is preserved there, while row 29 is converted
# to the numeric phone value (9876543) with points reset to 0.00.
$ws.Rows.Item(30).Insert()

# Row 30: keep the original text form of the phone number, points = 0.00
$ws.Cells.Item(30, 1).NumberFormat = "@"
$ws.Cells.Item(30, 1).Value = "09876543"
$ws.Cells.Item(30, 3).Value = 0

# Row 29: phone number becomes a true numeric value (no leading zero), points = 0.00
$ws.Cells.Item(29, 1).Value = 9876543
$ws.Cells.Item(29, 3).Value = 0
